$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("IntelliJ")

# Add the new row of material (Eclipselink-related IntelliJ shortcut)
$ws.Range("A10").Value = "Show methods in a call in a popup (like ctrl + o in eclipse) "
$ws.Range("B10").Value = "Ctrl + F12"

# Move the active selection to the next empty row, like Excel would after
# the user finishes entering data on row 10
$ws.Range("A11").Select()
